$d = $word.ActiveDocument

# Locate the "Ver no Jupiter ..." paragraph and the "(c) 2020 ... Contact: ..." paragraph
# that follows it near the end of the document (the site-footer boilerplate that was
# removed from the page when it was rebuilt).
$jupiterIndex = -1
$copyrightIndex = -1

for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $t = $d.Paragraphs.Item($i).Range.Text
    if ($jupiterIndex -eq -1 -and $t -like "*Ver no Jupiter*") {
        $jupiterIndex = $i
    }
    if ($copyrightIndex -eq -1 -and $t -like "*Contact: luizeleno@usp.br*") {
        $copyrightIndex = $i
    }
}

if ($jupiterIndex -ne -1 -and $copyrightIndex -ne -1) {
    # Also remove the blank paragraph that sits right before "Ver no Jupiter ...",
    # so the spacing around the remaining text stays the same as before this block
    # was ever inserted.
    $startIndex = $jupiterIndex
    if ($jupiterIndex -gt 1) {
        $prevText = $d.Paragraphs.Item($jupiterIndex - 1).Range.Text
        if ($prevText.Trim().Length -eq 0) {
            $startIndex = $jupiterIndex - 1
        }
    }

    $startPara = $d.Paragraphs.Item($startIndex)
    $endPara = $d.Paragraphs.Item($copyrightIndex)

    $r = $d.Range($startPara.Range.Start, $endPara.Range.End)
    $r.Delete()
}
